$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.048758552317011095
$ws.Range("B1").Value = 0.048475628702846052
$ws.Range("A2").Value = 0.04378355625694752
$ws.Range("B2").Value = -0.044483208807333696
$ws.Range("A3").Value = 0.14741164639237425
$ws.Range("B3").Value = -0.14811064177163402
$ws.Range("A4").Value = -0.13587566805062323
$ws.Range("B4").Value = 0.13526426707265671
$ws.Range("A5").Value = -0.12926426774933919
$ws.Range("B5").Value = 0.12805232404317746
$ws.Range("A6").Value = -0.02722590228885613
$ws.Range("B6").Value = 0.02721344554009697
$ws.Range("A7").Value = -0.034444881079837231
$ws.Range("B7").Value = 0.034379567614827522
$ws.Range("A8").Value = -0.014379568433783518
$ws.Range("B8").Value = 0.014363368977620361
$ws.Range("A9").Value = -0.0083633696793903312
$ws.Range("B9").Value = 0.008357607113737231
$ws.Range("A10").Value = -0.0023576078170535197
$ws.Range("B10").Value = 0.0023590319099682233
$ws.Range("A11").Value = 0.0021409673994021716
$ws.Range("B11").Value = -0.0021402084691608536
$ws.Range("A12").Value = -0.069955880620222732
$ws.Range("B12").Value = 0.069380562235142218
$ws.Range("A13").Value = -0.063380562956345088
$ws.Range("B13").Value = 0.063205210815896251
$ws.Range("A14").Value = -0.051205211594115951
$ws.Range("B14").Value = 0.051084393703285436
$ws.Range("A15").Value = -0.045084394435043862
$ws.Range("B15").Value = 0.044962304944857756
$ws.Range("A16").Value = -0.038962305681870202
$ws.Range("B16").Value = 0.038800221320141492
$ws.Range("A17").Value = -0.032800222066000195
$ws.Range("B17").Value = 0.032733907794731287
$ws.Range("A18").Value = -0.082259185672374002
$ws.Range("B18").Value = 0.082142453597722209
$ws.Range("A19").Value = -0.073142454303694482
$ws.Range("B19").Value = 0.072210797725117182
$ws.Range("A20").Value = -0.06321079844969546
$ws.Range("B20").Value = 0.06300399016671232
$ws.Range("A21").Value = -0.054003990895116427
$ws.Range("B21").Value = 0.053714722284969429
$ws.Range("A22").Value = -0.093928438286654981
$ws.Range("B22").Value = 0.093622077988404939
$ws.Range("A23").Value = -0.084622078696122927
$ws.Range("B23").Value = 0.084124201636380569
$ws.Range("A24").Value = -0.042124202638138541
$ws.Range("B24").Value = 0.041999998992919174
$ws.Range("A25").Value = -0.094826422449084902
$ws.Range("B25").Value = 0.094592131748349573
$ws.Range("A26").Value = -0.088592132445818095
$ws.Range("B26").Value = 0.088288779152815522
$ws.Range("A27").Value = -0.082288779853921135
$ws.Range("B27").Value = 0.081245512675295117
$ws.Range("A28").Value = -0.075245513391860364
$ws.Range("B28").Value = 0.074522829880336872
$ws.Range("A29").Value = -0.062522830658368278
$ws.Range("B29").Value = 0.062165253317695957
$ws.Range("A30").Value = -0.042165254168862631
$ws.Range("B30").Value = 0.042018474904665482
$ws.Range("A31").Value = -0.027018475720170798
$ws.Range("B31").Value = 0.02700054867840862
$ws.Range("A32").Value = -0.0060005495464388048
$ws.Range("B32").Value = 0.0059999992600223706

$ws.Columns.Item(2).ColumnWidth = 14.6
